$d = $word.ActiveDocument

# The heading paragraph (2nd paragraph) currently holds the "_GoBack" bookmark
# at its end (after the two split runs). We want it moved to the very start
# of the paragraph, before the (to-be-merged) run. Re-adding a bookmark under
# an already-existing name only relocates its start marker reliably while the
# document is still in its original shape, so: delete the old one, then add
# the new one at the paragraph's start *before* doing any further edits.
$d.Bookmarks("_GoBack").Delete()

$heading = $d.Paragraphs(2)
$startOfHeading = $d.Range($heading.Range.Start, $heading.Range.Start)
$d.Bookmarks.Add("_GoBack", $startOfHeading)

# Merge the two split bold runs ("COMMITING & PUSHING CHANGES TO " + "THE
# CSX415-assignments REPO") back into a single contiguous run/text node by
# replacing the combined text with itself.
$d.Content.Find.Execute(
    "COMMITING & PUSHING CHANGES TO THE CSX415-assignments REPO",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "COMMITING & PUSHING CHANGES TO THE CSX415-assignments REPO", 2)

# Remove the leading empty paragraph.
$d.Paragraphs(1).Range.Delete()
